$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "b.md" has now been handed off for localization. Update the Overview sheet
# and both the zh-cn / de-de localization-status sheets accordingly.
# ---------------------------------------------------------------------------

$handoffStatus   = "Ready for handoff"
$handoffDateTime = "2016-08-26 22:36:44"

# --- Overview sheet: row 3 corresponds to b.md ---------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E3").Value = $handoffStatus
$overview.Range("F3").Value = $handoffStatus
$overview.Range("G3").Value = $handoffDateTime

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")

# Row 2 = a.md : status simply changes to "Ready for handoff"
$zhcn.Range("C2").Value = $handoffStatus

# Row 3 = b.md : new handoff file generated for localization
$zhcn.Range("C3").Value = $handoffStatus
# Leading apostrophe forces text storage so "False" isn't auto-coerced to a Boolean
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-26 22:36:39"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3205c7784bfbe263d885520fc8b0477cd564cb79/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0446139d6f389bbea27781132fa0267257940e97/e2e/b.md."

# Column P widened to fit the new error message
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet -------------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")

# Row 2 = a.md : status simply changes to "Ready for handoff"
$dede.Range("C2").Value = $handoffStatus

# Row 3 = b.md : new handoff file generated for localization
$dede.Range("C3").Value = $handoffStatus
# Leading apostrophe forces text storage so "False" isn't auto-coerced to a Boolean
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $handoffDateTime
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3205c7784bfbe263d885520fc8b0477cd564cb79/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0446139d6f389bbea27781132fa0267257940e97/e2e/b.md."

# Column P widened to fit the new error message
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
